$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# Section header (row 23)
$ws.Range("A23").Value = "aaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaa"

# Sub header (row 25)
$ws.Range("A25").Value = "baza gtzan, broj žanrova 10 (svi), značajka mfcc, klasifikator kresvm(novi), treniranje 80%"

# Table header (row 27)
$ws.Range("B27").Value = "C"
$ws.Range("C27").Value = "gama/poli"
$ws.Range("D27").Value = "rezultat"
$ws.Range("A27").Value = "vrsta kernela"

# Data rows 28-42
$cValues = @(0.1, 0.5, 0.6, 0.7, 0.8, 0.9, 1, 1.1, 1.2, 1.3, 1.5, 2, 3, 4, 5)
$resValues = @(0.31, 0.475, 0.48, 0.47, 0.445, 0.43, 0.425, 0.435, 0.43, 0.445, 0.42, 0.405, 0.395, 0.375, 0.395)

for ($i = 0; $i -lt $cValues.Length; $i++) {
    $row = 28 + $i
    $ws.Range("A$row").Value = "linearni (1)"
    $ws.Range("B$row").Value = $cValues[$i]
    $ws.Range("D$row").Value = $resValues[$i]
    $ws.Range("D$row").Style = "Percent"
    $ws.Range("D$row").NumberFormat = "0.00%"
}

# Column widths
$ws.Columns.Item(3).ColumnWidth = 13.28515625
$ws.Columns.Item(4).ColumnWidth = 9.42578125

# Update view: scroll to A23, select D32
$ws.Range("D32").Select()
$excel.ActiveWindow.ScrollRow = 23
